# Generate Report for Handoff
# Updates status/date/priority/error-detail values across the Overview,
# zh-cn and de-de sheets to reflect the new handoff report.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

$overviewDateNew = "2017-02-22 08:44:07"
$zhcnHandoffDateNew = "2017-02-22 08:43:50"
$dedeHandoffDateNew = "2017-02-22 08:44:07"

$priorityNew = "mt"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/0344edee42689cb939822776aba7ffc2ad7e9588/e2e/69777f57-745b-46c0-93a6-146988fa487e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/8b56ee1066417d60bd8488d6895fc56cc228d33a/e2e/69777f57-745b-46c0-93a6-146988fa487e.md."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("G2").Value = $overviewDateNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $overviewDateNew

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("E2").Value = $priorityNew
$wsZhCn.Range("H2").Value = $zhcnHandoffDateNew
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("E3").Value = $priorityNew
$wsZhCn.Range("H3").Value = $zhcnHandoffDateNew
$wsZhCn.Range("R3").Value = $errorDetail

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("E2").Value = $priorityNew
$wsDeDe.Range("H2").Value = $dedeHandoffDateNew
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("E3").Value = $priorityNew
$wsDeDe.Range("H3").Value = $dedeHandoffDateNew
$wsDeDe.Range("R3").Value = $errorDetail

# --- Column width adjustments (Status columns narrower, Error Detail wider) ---
# Excel quantizes ColumnWidth to whole-pixel boundaries on write, so these are
# the closest attainable values to the authored widths (status cols ~17.2,
# error-detail cols = 40).
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25

$wsZhCn.Columns.Item(3).ColumnWidth = 16.25
$wsZhCn.Columns.Item(18).ColumnWidth = 39.1

$wsDeDe.Columns.Item(3).ColumnWidth = 16.25
$wsDeDe.Columns.Item(18).ColumnWidth = 39.1
